$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 previously held the shared string "R40"; the row now records the
# rule's sequence number "1" (stored as text, matching the column's
# existing shared-string cell type).
$ws.Range("B11").Value = "1"
